# Fixed update to excel issue
# - Rename "Requested quantity" header on Weekly/Monthly sheets
# - Add a new "PO Forecast" sheet with forecast data

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# Rename the "Requested quantity" headers.
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Header row.
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Match the header formatting used on the other sheets (bold, centered, bordered).
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows.
$newSheet.Range("A2").Value = 45004.99999999999
$newSheet.Range("B2").Value = 9
$newSheet.Range("C2").Value = 3.379353855417127
$newSheet.Range("D2").Value = 13.97391323031135
$newSheet.Range("A3").Value = 45011.99999999999
$newSheet.Range("B3").Value = 7
$newSheet.Range("C3").Value = 2.120951286721005
$newSheet.Range("D3").Value = 12.69544690143536
$newSheet.Range("A4").Value = 45018.99999999999
$newSheet.Range("B4").Value = 6
$newSheet.Range("C4").Value = 0.6525191882383053
$newSheet.Range("D4").Value = 11.08155888529303
$newSheet.Range("A5").Value = 45025.99999999999
$newSheet.Range("B5").Value = 4
$newSheet.Range("C5").Value = -0.9697211782020976
$newSheet.Range("D5").Value = 9.57023254631372
$newSheet.Range("A6").Value = 45032.99999999999
$newSheet.Range("B6").Value = 3
$newSheet.Range("C6").Value = -2.407126209035489
$newSheet.Range("D6").Value = 8.326931211088716
$newSheet.Range("A7").Value = 45039.99999999999
$newSheet.Range("B7").Value = 1
$newSheet.Range("C7").Value = -3.685210434069158
$newSheet.Range("D7").Value = 6.615634021843968
$newSheet.Range("A8").Value = 45046.99999999999
$newSheet.Range("B8").Value = 0
$newSheet.Range("C8").Value = -5.689607742696194
$newSheet.Range("D8").Value = 5.110471888939085
$newSheet.Range("A9").Value = 45053.99999999999
$newSheet.Range("B9").Value = 0
$newSheet.Range("C9").Value = -6.840256007512473
$newSheet.Range("D9").Value = 3.973677401675564
$newSheet.Range("A10").Value = 45060.99999999999
$newSheet.Range("B10").Value = 0
$newSheet.Range("C10").Value = -8.265989246474241
$newSheet.Range("D10").Value = 2.273801138360441
$newSheet.Range("A11").Value = 45067.99999999999
$newSheet.Range("B11").Value = 0
$newSheet.Range("C11").Value = -9.656825782324368
$newSheet.Range("D11").Value = 0.4548692515485713
$newSheet.Range("A12").Value = 45074.99999999999
$newSheet.Range("B12").Value = 0
$newSheet.Range("C12").Value = -11.05450305012978
$newSheet.Range("D12").Value = -1.162696260038573

# Match the date formatting used in column A of the other sheets.
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A12").PasteSpecial(-4122)

$newSheet.Range("A1").Select()
